$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 115; this shifts all rows 115..191 down to 116..192
# and automatically grows the sheet's used range/dimension from R191 to R192.
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row 115 with the new record.
# Columns A,B,C,E,F,G,H,I,N,Q,R are the same constant values used by every
# other data row in this sheet; D,J,K,L,M,O,P hold this record's specific data.
$ws.Cells.Item(115, 1).Value = 6
$ws.Cells.Item(115, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(115, 3).Value = "Metropolitana"
$ws.Cells.Item(115, 4).Value = 44673
$ws.Cells.Item(115, 5).Value = 13
$ws.Cells.Item(115, 6).Value = 100112029
$ws.Cells.Item(115, 7).Value = "Orégano"
$ws.Cells.Item(115, 8).Value = "Sin especificar"
$ws.Cells.Item(115, 9).Value = "Primera"
$ws.Cells.Item(115, 10).Value = 33
$ws.Cells.Item(115, 11).Value = 14000
$ws.Cells.Item(115, 12).Value = 15000
$ws.Cells.Item(115, 13).Value = 14455
$ws.Cells.Item(115, 14).Value = "$/docena de atados"
$ws.Cells.Item(115, 15).Value = "Región Metropolitana"
$ws.Cells.Item(115, 16).Value = 4818
$ws.Cells.Item(115, 17).Value = 3
$ws.Cells.Item(115, 18).Value = "Hortaliza"
